# Auto-generated edit script: update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "88.692.04"
$ws.Range("E2").Value = "  +8.94%  "

$ws.Range("D3").Value = "3.341.97"
$ws.Range("E3").Value = "  +5.10%  "

$ws.Range("E4").Value = "  +0.10%  "

$c = $ws.Range("D5")
$c.Value = "'219.74"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +5.63%  "

$c = $ws.Range("D6")
$c.Value = "'651.71"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.00%  "

$c = $ws.Range("D7")
$c.Value = "'0.393"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +32.73%  "

$c = $ws.Range("D8")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.00%  "

$c = $ws.Range("D9")
$c.Value = "'0.607"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.66%  "

$ws.Range("D10").Value = "3.336.25"
$ws.Range("E10").Value = "  +5.00%  "

$c = $ws.Range("D11")
$c.Value = "'0.590"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.41%  "

$c = $ws.Range("D12")
$c.Value = "'0.0000280"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +6.71%  "

$ws.Range("E13").Value = "  +1.63%  "

$c = $ws.Range("D14")
$c.Value = "'35.25"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +10.41%  "

$ws.Range("D15").Value = "3.956.18"
$ws.Range("E15").Value = "  +5.13%  "

$c = $ws.Range("D16")
$c.Value = "'5.48"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.86%  "

$ws.Range("D17").Value = "88.487.95"
$ws.Range("E17").Value = "  +8.65%  "

$ws.Range("D18").Value = "3.324.78"
$ws.Range("E18").Value = "  +4.51%  "

$c = $ws.Range("D19")
$c.Value = "'14.66"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.37%  "

$c = $ws.Range("D20")
$c.Value = "'3.16"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.60%  "

$c = $ws.Range("D21")
$c.Value = "'9.77"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +5.68%  "

$c = $ws.Range("D22")
$c.Value = "'457.52"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +4.21%  "

$c = $ws.Range("D23")
$c.Value = "'5.51"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +6.06%  "

$c = $ws.Range("D24")
$c.Value = "'7.41"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +3.83%  "

$c = $ws.Range("D25")
$c.Value = "'5.60"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +6.42%  "

$c = $ws.Range("D26")
$c.Value = "'12.82"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +14.06%  "

$ws.Range("D27").Value = "3.515.47"
$ws.Range("E27").Value = "  +5.13%  "

$c = $ws.Range("D28")
$c.Value = "'78.60"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.21%  "

$ws.Range("E29").Value = "  +3.02%  "

$c = $ws.Range("D30")
$c.Value = "'0.200"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +41.39%  "

$ws.Range("E31").Value = "  +0.14%  "

$c = $ws.Range("D32")
$c.Value = "'9.41"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.69%  "

$c = $ws.Range("D33")
$c.Value = "'595.28"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +6.26%  "

$ws.Range("E34").Value = "  +5.94%  "

$c = $ws.Range("D35")
$c.Value = "'0.998"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.29%  "

$c = $ws.Range("D36")
$c.Value = "'2.14"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +4.82%  "

$c = $ws.Range("D37")
$c.Value = "'7.27"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +21.24%  "

$c = $ws.Range("D38")
$c.Value = "'0.146"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -4.56%  "

$c = $ws.Range("D39")
$c.Value = "'23.39"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.11%  "

$c = $ws.Range("D40")
$c.Value = "'2.17"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +6.54%  "

$c = $ws.Range("D41")
$c.Value = "'0.422"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.54%  "

$c = $ws.Range("D42")
$c.Value = "'21.90"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +5.54%  "

$c = $ws.Range("D43")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("E44").Value = "  -0.31%  "

$c = $ws.Range("D45")
$c.Value = "'159.03"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.86%  "

$ws.Range("E46").Value = "  +0.02%  "

$c = $ws.Range("D47")
$c.Value = "'190.70"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.32%  "

$ws.Range("E48").Value = "  +6.69%  "

$ws.Range("E49").Value = "  +4.93%  "

$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D50")
$c.Value = "'4.43"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.33%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D51")
$c.Value = "'0.785"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.09%  "

